# Applies the "Add function to see all bridges that ever existed" edit
# to the final-project-proposal.docx Functional Specifications /
# Anticipated User Interaction sections.

$d = $word.ActiveDocument

function Find-ParagraphByText($text) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find paragraph containing: $text"
    }
    return $rng.Paragraphs(1)
}

# ---------------------------------------------------------------------
# 1. Insert a new bullet "The city's current outline will be displayed
#    for reference." right after "A legend will inform users of the
#    bridge types." (i.e. right before the "Functional Specifications"
#    heading).
# ---------------------------------------------------------------------
$legendPara = Find-ParagraphByText("A legend will inform users of the bridge types.")
$legendPara.Range.InsertParagraphAfter()
$newPara = $legendPara.Next()
$newPara.Range.Text = "The city's current outline will be displayed for reference."

# ---------------------------------------------------------------------
# 2. "A search box will allow users to search for a specific bridge" ->
#    add " in the selected year" right before the trailing ". A possible
#    development example..." sentence.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "A search box will allow users to search for a specific bridge. A possible development example",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A search box will allow users to search for a specific bridge in the selected year. A possible development example",
    2) | Out-Null

# ---------------------------------------------------------------------
# 3. Remove the whole "second search box / geocoding" bullet paragraph
#    entirely (including its hyperlink).
# ---------------------------------------------------------------------
$geoPara = Find-ParagraphByText("If possible without navigating users outside of the map bounds")
$geoPara.Range.Delete()

# ---------------------------------------------------------------------
# 4. Insert a new bullet "Bridge points will appear lighter the older
#    they become." right before "Users will be able to click on a
#    bridge to get more information..." paragraph.
# ---------------------------------------------------------------------
$clickPara = Find-ParagraphByText("Users will be able to click on a bridge to get more information")
$clickPara.Range.InsertParagraphBefore()
# NOTE: InsertParagraphBefore() keeps $clickPara anchored to the (now
# empty) newly-created paragraph; the original content shifts to
# $clickPara.Next().
$clickPara.Range.Text = "Bridge points will appear lighter the older they become."

# ---------------------------------------------------------------------
# 5. Insert a new bullet about the all-bridges list feature right after
#    "Users will be able to hover over a bridge to get its name, year
#    built, and year demolished." and add the marker-list-click link.
# ---------------------------------------------------------------------
$hoverPara = Find-ParagraphByText("Users will be able to hover over a bridge to get its name, year built, and year demolished.")
$hoverPara.Range.InsertParagraphAfter()
$listPara = $hoverPara.Next()
$listPara.Range.Text = "If possible, bridges from all years will be displayed in a list. If a user clicks a bridge from the list, the map will pan to it, replace the list with the detail panel for that bridge, and zoom the slider to the year the bridge was built. If the user closes the detail panel, the list will reappear in its place. See LINKPLACEHOLDER. "

$linkRange = $listPara.Range.Duplicate()
$linkRange.Find.Execute("LINKPLACEHOLDER", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Hyperlinks.Add($linkRange, "https://www.mapbox.com/mapbox.js/example/v1.0.0/marker-list-click/", [Type]::Missing, [Type]::Missing, "https://www.mapbox.com/mapbox.js/example/v1.0.0/marker-list-click/") | Out-Null

# ---------------------------------------------------------------------
# 6. Remove "Search for a specific address or location within the city"
#    bullet entirely (Anticipated User Interaction section).
# ---------------------------------------------------------------------
$searchAddrPara = Find-ParagraphByText("Search for a specific address or location within the city")
$searchAddrPara.Range.Delete()

# ---------------------------------------------------------------------
# 7. "Click a bridge to see a photo or drawing and description" ->
#    "Click a bridge on the map or list of all bridges to see a photo
#    or drawing and description"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Click a bridge to see a photo or drawing and description",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Click a bridge on the map or list of all bridges to see a photo or drawing and description",
    2) | Out-Null

Write-Output "done"
